$d = $word.ActiveDocument

# 1. Title: "Global Variables" -> "Global " + "Object" as two separate runs.
#    A plain Find/Replace would coalesce the result into a single run, so we
#    locate the word "Variables" and overwrite just that range with a
#    brand-new run ("Object") via raw OOXML, which keeps it distinct from the
#    pre-existing "Global " run instead of merging into it.
$titleRange = $d.Content
$titleRange.Find.Execute("Variables", $true, $true, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$titleRange.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Object</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# 2. Table header cell: "Variable" -> "Object property"
$d.Content.Find.Execute("Variable", $true, $true, $false, $false, $false, $true, 1, $false, "Object property", 2) | Out-Null

# 3. Merge the two runs "workT" + "ime" into a single run "workTime"
$d.Content.Find.Execute("workTime", $true, $true, $false, $false, $false, $true, 1, $false, "workTime", 2) | Out-Null

# 4. Merge the two runs "break" + "TimeLeft" into a single run "breakTimeLeft"
$d.Content.Find.Execute("breakTimeLeft", $true, $true, $false, $false, $false, $true, 1, $false, "breakTimeLeft", 2) | Out-Null

# 5. Add a new row "Play" / "boolean" / "True when playing" at the end of the table
$t = $d.Tables(1)
$newRow = $t.Rows.Add()
$idx = $t.Rows.Count

$cell1 = $t.Cell($idx, 1)
$cell1.Range.Text = "Play"

$cell3 = $t.Cell($idx, 3)
$cell3.Range.Text = "True when playing"

$cell2 = $t.Cell($idx, 2)
$cell2.Range.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>boolean</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null
